# Applies the "Added BWP Bootstrap files" update to VT-Data-Prod.xlsx.
# This refreshes the recorded Prod/Demo run-timestamps (columns B = DateProd,
# D = DateDemo) on each of the 9 bootstrap result sheets with the latest
# Katalon execution timestamps, exactly as captured by the source commit.

$wb = $excel.ActiveWorkbook

function Set-Cell($ws, [string]$addr, [string]$val) {
    $ws.Range($addr).Value = $val
}

# --- VT-SaleVoid-DualCF-Generic ---
$ws = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
Set-Cell $ws "B2" "Mon Nov 10 17:53:02 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 17:53:51 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 17:54:42 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 17:55:37 IST 2025"

# --- VT-SaleVoid-NoCF-Generic ---
$ws = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
Set-Cell $ws "B2" "Mon Nov 10 17:56:33 IST 2025"
Set-Cell $ws "D2" "Thu Nov 06 23:49:33 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 17:57:31 IST 2025"
Set-Cell $ws "D3" "Thu Nov 06 23:50:44 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 17:58:26 IST 2025"
Set-Cell $ws "D4" "Thu Nov 06 23:51:33 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 17:59:46 IST 2025"
Set-Cell $ws "D5" "Thu Nov 06 23:52:20 IST 2025"

# --- VT-SaleVoid-SingleCF-Generic ---
$ws = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
Set-Cell $ws "B2" "Mon Nov 10 18:00:53 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 18:01:44 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 18:02:50 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 18:03:56 IST 2025"

# --- VT-SaleCredit-DualCF-Generic ---
$ws = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
Set-Cell $ws "B2" "Mon Nov 10 17:45:01 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 17:46:12 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 17:47:00 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 17:48:01 IST 2025"

# --- VT-SaleCredit-NoCF-Generic ---
$ws = $wb.Worksheets.Item("VT-SaleCredit-NoCF-Generic")
Set-Cell $ws "B2" "Thu Nov 06 23:36:54 IST 2025"
Set-Cell $ws "B3" "Thu Nov 06 23:37:42 IST 2025"
Set-Cell $ws "B4" "Thu Nov 06 23:38:31 IST 2025"
Set-Cell $ws "B5" "Thu Nov 06 23:39:20 IST 2025"

# --- VT-SaleCredit-SingleCF-Generic ---
$ws = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
Set-Cell $ws "B2" "Mon Nov 10 17:49:16 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 17:50:31 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 17:51:16 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 17:52:07 IST 2025"

# --- VT-AuthCapCredit-Generic ---
$ws = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
Set-Cell $ws "B2" "Mon Nov 10 16:45:24 IST 2025"
Set-Cell $ws "D2" "Thu Nov 06 22:45:08 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 16:46:41 IST 2025"
Set-Cell $ws "D3" "Thu Nov 06 22:46:18 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 16:47:50 IST 2025"
Set-Cell $ws "D4" "Thu Nov 06 22:47:29 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 16:49:13 IST 2025"
Set-Cell $ws "D5" "Thu Nov 06 22:48:49 IST 2025"
Set-Cell $ws "B6" "Mon Nov 10 16:50:47 IST 2025"
Set-Cell $ws "D6" "Thu Nov 06 22:50:14 IST 2025"
Set-Cell $ws "B7" "Mon Nov 10 16:52:13 IST 2025"
Set-Cell $ws "D7" "Thu Nov 06 22:51:21 IST 2025"

# --- VT-AuthCapVoid-Generic ---
$ws = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
Set-Cell $ws "B2" "Mon Nov 10 17:01:03 IST 2025"
Set-Cell $ws "D2" "Fri Nov 07 19:23:29 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 17:02:20 IST 2025"
Set-Cell $ws "D3" "Fri Nov 07 19:25:06 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 17:03:41 IST 2025"
Set-Cell $ws "D4" "Fri Nov 07 19:26:30 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 17:04:50 IST 2025"
Set-Cell $ws "D5" "Fri Nov 07 19:27:41 IST 2025"
Set-Cell $ws "B6" "Mon Nov 10 17:06:00 IST 2025"
Set-Cell $ws "D6" "Fri Nov 07 19:29:01 IST 2025"
Set-Cell $ws "B7" "Mon Nov 10 17:07:21 IST 2025"
Set-Cell $ws "D7" "Fri Nov 07 19:30:23 IST 2025"

# --- VT-ManualAuthCapture-Generic ---
$ws = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")
Set-Cell $ws "B2" "Mon Nov 10 17:38:52 IST 2025"
Set-Cell $ws "B3" "Mon Nov 10 17:39:46 IST 2025"
Set-Cell $ws "B4" "Mon Nov 10 17:40:51 IST 2025"
Set-Cell $ws "B5" "Mon Nov 10 17:42:26 IST 2025"
Set-Cell $ws "B6" "Mon Nov 10 17:43:18 IST 2025"
Set-Cell $ws "B7" "Mon Nov 10 17:44:08 IST 2025"
